# Update from MV -datos-
# Updates row 74 (2021Q1 "01-01-2021") with revised figures and appends a
# new row 75 for the "01-04-2021" (2021Q2) period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 74 values that changed ---
$ws.Range("B74").Value = -1857
$ws.Range("E74").Value = -63
$ws.Range("F74").Value = 2035
$ws.Range("L74").Value = -3485
$ws.Range("M74").Value = -415
$ws.Range("Q74").Value = 958
$ws.Range("S74").Value = 1243
$ws.Range("U74").Value = -5
$ws.Range("V74").Value = 3364
$ws.Range("W74").Value = 1896
$ws.Range("X74").Value = 1694
$ws.Range("Z74").Value = -229
$ws.Range("AA74").Value = -1096
$ws.Range("AB74").Value = -1988
$ws.Range("AC74").Value = 713
$ws.Range("AE74").Value = 98

# --- Add new row 75 ("01-04-2021") ---
# A75 holds the text label "01-04-2021". Assigning it directly would be
# auto-recognized as a date literal by Excel's smart input, so instead we
# compute it as a text formula result and then paste-special just the
# value back over itself; this yields a plain shared-string text cell
# identical to how the other period labels (A2:A74) are stored.
$ws.Range("A75").Formula = "=""01-04-2021"""
$ws.Range("A75").Copy() | Out-Null
$ws.Range("A75").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B75").Value = -8496
$ws.Range("C75").Value = -8638
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 142
$ws.Range("F75").Value = 2614
$ws.Range("G75").Value = 4772
$ws.Range("H75").Value = 4446
$ws.Range("I75").Value = 248
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 78
$ws.Range("L75").Value = 2930
$ws.Range("M75").Value = -248
$ws.Range("N75").Value = 829
$ws.Range("O75").Value = 102
$ws.Range("P75").Value = 2247
$ws.Range("Q75").Value = -6622
$ws.Range("R75").Value = 0
$ws.Range("S75").Value = -6609
$ws.Range("T75").Value = -121
$ws.Range("U75").Value = 108
$ws.Range("V75").Value = 1534
$ws.Range("W75").Value = 61
$ws.Range("X75").Value = 1474
$ws.Range("Y75").Value = -22
$ws.Range("Z75").Value = 21
$ws.Range("AA75").Value = 1087
$ws.Range("AB75").Value = -3249
$ws.Range("AC75").Value = 982
$ws.Range("AD75").Value = 116
$ws.Range("AE75").Value = 3238
